$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46010
$ws.Range("B2").Value = 89.47
$ws.Range("C2").Value = 86.56
$ws.Range("D2").Value = 83.31999999999999
$ws.Range("E2").Value = 81.40000000000001
$ws.Range("F2").Value = 83.17
$ws.Range("G2").Value = 88.42
$ws.Range("H2").Value = 99.16
$ws.Range("I2").Value = 113.2
$ws.Range("J2").Value = 113.6
$ws.Range("K2").Value = 109.06
$ws.Range("L2").Value = 105.11
$ws.Range("M2").Value = 104.51
$ws.Range("N2").Value = 103.43
$ws.Range("O2").Value = 104.3
$ws.Range("P2").Value = 107.21
$ws.Range("Q2").Value = 107.74
$ws.Range("R2").Value = 103.05
$ws.Range("S2").Value = 103.83
$ws.Range("T2").Value = 113.79
$ws.Range("U2").Value = 116.1
$ws.Range("V2").Value = 117.51
$ws.Range("W2").Value = 115.55
$ws.Range("X2").Value = 110.88
$ws.Range("Y2").Value = 105.9
$ws.Range("Z2").Value = 102.76
$ws.Range("AB2").Value = 112.46
$ws.Range("AD2").Value = 116.53
$ws.Range("AF2").Value = 114.94
$ws.Range("AG2").Value = "0h-6h"
